$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Add the new shared strings in the order matching the target file:
# 37 -> Google Sheet Key (K3), 38 -> OCT Credentials path (I3), 39 -> Output file path (C3)
$ws.Range("K3").Value = "19MGCdQZPN6ucU9prUqnhfxYmlGt2QO20w_IT6w3pVXA"
$ws.Range("I3").Value = "F:\Tools\Python\credentials\wp_oct_api\custom-healer-187616-2e4a86fd3dc2.json"

# C3 needs to carry the same format/style as C2 (the yellow highlighted input cell),
# so copy C2's formatting onto C3 before setting its value.
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "F:\Clients\SGBU\Tools\OCT\SGBU_OCT_output_22Mar24_conversion_test.csv"

# Update the active selection to C2, matching the saved workbook state.
$ws.Range("C2").Select()
